$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N,O,P -> O,P,Q)
$ws.Columns("N:N").Insert()

# The newly inserted column should inherit the width of the column to its
# left (column M), matching Excel's native "insert column" behaviour.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with a new selection.
$ws.Select()
$ws.Range("S7").Select()
